$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 10)
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = [double]"8.943567213078773E-15"
$ws.Range("E2").Value = [double]"8.943567213078773E-15"

# Row 3 (Control 5)
$ws.Range("D3").Value = [double]"0.9999997511863719"
$ws.Range("E3").Value = [double]"0.9999997511863719"

# Row 4 (MDD 41)
$ws.Range("D4").Value = [double]"5.00661137352761E-09"
$ws.Range("E4").Value = [double]"0.9999999949933887"

# Row 5 (MDD 8)
$ws.Range("D5").Value = [double]"0.999999991577737"
$ws.Range("E5").Value = [double]"8.422262998131202E-09"

# Row 6 (MDD 15)
$ws.Range("D6").Value = [double]"8.840944593218144E-14"
$ws.Range("E6").Value = [double]"0.9999999999999116"

# Row 7 (MDD 16)
$ws.Range("D7").Value = [double]"0.9995679070259491"
$ws.Range("E7").Value = [double]"0.0004320929740508905"

# Row 8 (MDD 33)
$ws.Range("D8").Value = [double]"0.9999807855388351"
$ws.Range("E8").Value = [double]"1.92144611649292E-05"
$ws.Range("F8").Value = [double]"9.196616172790527"
$ws.Range("G8").Value = [double]"0.5714285714285714"
